# feat: add 2022-Q3 data
#
# Insert a new "2022-Q3" worksheet between "总计" and "2021-Q1", fill it
# with the quarter's holdings, and update the "总计" summary sheet so its
# existing "2021-Q1" row becomes row 3 while row 2 now reports 2022-Q3.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Item("2021-Q1")

# Insert the new sheet right before "2021-Q1" so the tab order becomes
# 总计, 2022-Q3, 2021-Q1.
$wsQ3 = $wb.Worksheets.Add($wsQ1)
$wsQ3.Name = "2022-Q3"

# --- "总计" sheet: shift the existing 2021-Q1 row down to row 3, then
# overwrite row 2 with the new 2022-Q3 summary figures. ---
$wsTotal.Cells.Item(2,1).Copy($wsTotal.Cells.Item(3,1))
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(3,2).Value = "2021-Q1"
$wsTotal.Cells.Item(3,3).Value = 2
$wsTotal.Cells.Item(3,4).Value = 0.03

$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 1
$wsTotal.Cells.Item(2,4).Value = 0

# --- "2022-Q3" sheet: same layout/style as the other quarterly sheets
# (header row + A-column style borrowed from "总计"). ---
$wsTotal.Cells.Item(1,2).Copy($wsQ3.Range("B1:H1"))
$wsTotal.Cells.Item(2,1).Copy($wsQ3.Cells.Item(2,1))

$wsQ3.Cells.Item(1,2).Value = "基金代码"
$wsQ3.Cells.Item(1,3).Value = "基金名称"
$wsQ3.Cells.Item(1,4).Value = "基金规模"
$wsQ3.Cells.Item(1,5).Value = "股票总仓位"
$wsQ3.Cells.Item(1,6).Value = "仓位占比"
$wsQ3.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsQ3.Cells.Item(1,8).Value = "仓位排名"

$wsQ3.Cells.Item(2,1).Value = 0

# Columns B:G hold text-like figures ("005702", "0.20", …) that must stay
# text rather than being auto-coerced to numbers.
$dataText = $wsQ3.Range("B2:G2")
$dataText.NumberFormat = "@"
$wsQ3.Cells.Item(2,2).Value = "005702"
$wsQ3.Cells.Item(2,3).Value = "恒生前海港股通高股息低波动指数"
$wsQ3.Cells.Item(2,4).Value = "0.20"
$wsQ3.Cells.Item(2,5).Value = "94.22"
$wsQ3.Cells.Item(2,6).Value = "2.25"
$wsQ3.Cells.Item(2,7).Value = "0.0045"
$dataText.Style = "Normal"

$wsQ3.Cells.Item(2,8).Value = 10

# Restore "2021-Q1" as the active tab/selection (unaffected by this edit).
$wb.Worksheets.Item("2021-Q1").Activate()
